$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.196.54'
$ws.Range('D3').Value = '1.659.38'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = "'219.04"
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = "'0.5247"
$ws.Range('E6').Value = '  -1.03%  '
$ws.Range('D7').Value = "'1.007"
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('D8').Value = "'0.2630"
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'0.06303"
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('D10').Value = "'20.60"
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').Value = "'0.07809"
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').Value = "'4.486"
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').Value = '1.660.90'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '1.887.97'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = "'0.5540"
$ws.Range('D16').Value = '0.0₅8018'
$ws.Range('E16').Value = '  -1.61%  '
$ws.Range('D17').Value = "'65.16"
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '26.209.37'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').Value = "'4.626"
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').Value = "'196.23"
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('D22').Value = "'10.12"
$ws.Range('D23').Value = "'5.951"
$ws.Range('E23').Value = '  -1.19%  '
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').Value = "'145.55"
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = "'0.1202"
$ws.Range('E26').Value = '  -1.50%  '
$ws.Range('D27').Value = "'7.132"
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').Value = "'16.02"
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = "'0.05758"
$ws.Range('E30').Value = '  -2.08%  '
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').Value = "'3.481"
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('D33').Value = "'3.347"
$ws.Range('E33').Value = '  +2.64%  '
$ws.Range('D34').Value = "'1.580"
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('B35').Value = 'MXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D35').Value = "'2.806"
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'0.9522"
$ws.Range('E36').Value = '  -0.47%  '
$ws.Range('D37').Value = "'2.422"
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = "'0.5704"
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('D39').Value = "'0.01597"
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').Value = "'5.951"
$ws.Range('E40').Value = '  +1.86%  '
$ws.Range('D41').Value = '1.059.36'
$ws.Range('E41').Value = '  +1.50%  '
$ws.Range('D42').Value = "'0.8507"
$ws.Range('E42').Value = '  -1.47%  '
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').Value = "'103.11"
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = '1.799.25'
$ws.Range('D46').Value = "'58.09"
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('D47').Value = "'1.009"
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = "'0.4409"
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₈104'
$ws.Range('E49').Value = '  -1.83%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'8.015"
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.05206"
$ws.Range('E51').Value = '  +0.81%  '
